$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "'01/01/2023"
$ws.Range("C13").Value = "'01/01/2023"

$ws.Range("B15").Value = "3480026 - João Paulo Pascon"
$ws.Range("C15").Value = "3480026 - João Paulo Pascon"

$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
